$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D receive numeric-looking text (e.g. "211.05") that Excel
# would otherwise auto-convert to a number. Force text storage via the "@"
# number format, assign the value, then restore the default "Normal" style so
# the cell style index matches the original (no explicit style).
$dCells = @("D2","D3","D5","D6","D8","D9","D13","D14","D16","D17","D18","D19","D24","D25","D26","D27","D33","D34","D35","D37","D39","D41","D43","D44","D45","D46","D47","D48","D49","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.354.74'
$ws.Range("E2").Value = '  +2.40%  '

$ws.Range("D3").Value = '1.573.66'
$ws.Range("E3").Value = '  +0.52%  '

$ws.Range("E4").Value = '  +1.19%  '

$ws.Range("D5").Value = '211.05'
$ws.Range("E5").Value = '  +0.11%  '

$ws.Range("D6").Value = '0.490'
$ws.Range("E6").Value = '  -0.71%  '

$ws.Range("E7").Value = '  +1.29%  '

$ws.Range("D8").Value = '46.01'
$ws.Range("E8").Value = '  +4.39%  '

$ws.Range("D9").Value = '23.79'
$ws.Range("E9").Value = '  +3.33%  '

$ws.Range("E10").Value = '  -1.04%  '

$ws.Range("E11").Value = '  -0.62%  '

$ws.Range("E12").Value = '  +0.40%  '

$ws.Range("D13").Value = '1.799.16'
$ws.Range("E13").Value = '  +0.55%  '

$ws.Range("D14").Value = '1.585.22'
$ws.Range("E14").Value = '  +1.30%  '

$ws.Range("E15").Value = '  +0.83%  '

$ws.Range("D16").Value = '3.70'
$ws.Range("E16").Value = '  -0.84%  '

$ws.Range("D17").Value = '28.367.86'
$ws.Range("E17").Value = '  +2.44%  '

$ws.Range("D18").Value = '62.30'
$ws.Range("E18").Value = '  -1.62%  '

$ws.Range("D19").Value = '227.18'
$ws.Range("E19").Value = '  -0.63%  '

$ws.Range("E20").Value = '  -0.41%  '

$ws.Range("E21").Value = '  -1.48%  '

$ws.Range("E22").Value = '  +1.23%  '

$ws.Range("E23").Value = '  -4.16%  '

$ws.Range("D24").Value = '9.17'
$ws.Range("E24").Value = '  -1.57%  '

$ws.Range("D25").Value = '1.98'
$ws.Range("E25").Value = '  +3.58%  '

$ws.Range("D26").Value = '150.81'
$ws.Range("E26").Value = '  +0.39%  '

$ws.Range("D27").Value = '14.98'
$ws.Range("E27").Value = '  -1.48%  '

$ws.Range("E28").Value = '  -1.44%  '

$ws.Range("E29").Value = '  -2.08%  '

$ws.Range("E30").Value = '  +1.15%  '

$ws.Range("E31").Value = '  -1.40%  '

$ws.Range("E32").Value = '  -1.53%  '

$ws.Range("D33").Value = '3.20'
$ws.Range("E33").Value = '  -0.83%  '

$ws.Range("D34").Value = '3.12'
$ws.Range("E34").Value = '  +0.06%  '

$ws.Range("D35").Value = '1.391.95'
$ws.Range("E35").Value = '  -1.29%  '

$ws.Range("E36").Value = '  -1.18%  '

$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -2.87%  '

$ws.Range("E38").Value = '  +2.36%  '

$ws.Range("D39").Value = '2.58'
$ws.Range("E39").Value = '  +5.60%  '

$ws.Range("E40").Value = '  -1.32%  '

$ws.Range("D41").Value = '0.532'
$ws.Range("E41").Value = '  -1.79%  '

$ws.Range("E42").Value = '  +1.25%  '

$ws.Range("D43").Value = '0.794'
$ws.Range("E43").Value = '  -1.48%  '

$ws.Range("D44").Value = '5.61'
$ws.Range("E44").Value = '  -0.45%  '

$ws.Range("D45").Value = '1.86'
$ws.Range("E45").Value = '  +1.19%  '

$ws.Range("D46").Value = '0.979'
$ws.Range("E46").Value = '  +1.64%  '

$ws.Range("D47").Value = '62.36'
$ws.Range("E47").Value = '  -2.00%  '

$ws.Range("D48").Value = '1.710.03'
$ws.Range("E48").Value = '  +0.73%  '

$ws.Range("D49").Value = '85.78'
$ws.Range("E49").Value = '  -0.74%  '

$ws.Range("E50").Value = '  -0.64%  '

# Row 51: BabyDogeCoin -> BitcoinSV (name, link, price, volume all change)
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "38.97"
$ws.Range("E51").Value = '  -0.96%  '

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}